$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 237, shifting existing data (old rows 237-362) down to 239-364
$ws.Rows.Item(237).Resize(2).Insert()

# Fill in the two newly inserted rows with the new data
# Row 237 - Primera
$ws.Cells.Item(237, 1).Value = 9
$ws.Cells.Item(237, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(237, 3).Value = "Metropolitana"
$ws.Cells.Item(237, 4).Value = 44572
$ws.Cells.Item(237, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(237, 5).Value = 13
$ws.Cells.Item(237, 6).Value = 100114014
$ws.Cells.Item(237, 7).Value = "Betarraga"
$ws.Cells.Item(237, 8).Value = "Sin especificar"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 6100
$ws.Cells.Item(237, 11).Value = 90
$ws.Cells.Item(237, 12).Value = 100
$ws.Cells.Item(237, 13).Value = 95
$ws.Cells.Item(237, 14).Value = "`$/unidad"
$ws.Cells.Item(237, 15).Value = "Región Metropolitana"
$ws.Cells.Item(237, 16).Value = 95
$ws.Cells.Item(237, 17).Value = 1
$ws.Cells.Item(237, 18).Value = "Hortaliza"

# Row 238 - Segunda
$ws.Cells.Item(238, 1).Value = 9
$ws.Cells.Item(238, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(238, 3).Value = "Metropolitana"
$ws.Cells.Item(238, 4).Value = 44572
$ws.Cells.Item(238, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(238, 5).Value = 13
$ws.Cells.Item(238, 6).Value = 100114014
$ws.Cells.Item(238, 7).Value = "Betarraga"
$ws.Cells.Item(238, 8).Value = "Sin especificar"
$ws.Cells.Item(238, 9).Value = "Segunda"
$ws.Cells.Item(238, 10).Value = 2500
$ws.Cells.Item(238, 11).Value = 60
$ws.Cells.Item(238, 12).Value = 70
$ws.Cells.Item(238, 13).Value = 65
$ws.Cells.Item(238, 14).Value = "`$/unidad"
$ws.Cells.Item(238, 15).Value = "Región Metropolitana"
$ws.Cells.Item(238, 16).Value = 65
$ws.Cells.Item(238, 17).Value = 1
$ws.Cells.Item(238, 18).Value = "Hortaliza"
